$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "64.406.96"
Set-TextValue $ws.Range("E2") "  +3.10%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.077.69"
Set-TextValue $ws.Range("E3") "  +1.77%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.15%  "

# Row 5
Set-TextValue $ws.Range("D5") "559.24"
Set-TextValue $ws.Range("E5") "  +2.24%  "

# Row 6
Set-TextValue $ws.Range("D6") "145.67"
Set-TextValue $ws.Range("E6") "  +7.09%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  -0.10%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.075.51"
Set-TextValue $ws.Range("E8") "  +1.91%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.503"
Set-TextValue $ws.Range("E9") "  +1.60%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.153"
Set-TextValue $ws.Range("E10") "  +3.98%  "

# Row 11
Set-TextValue $ws.Range("D11") "6.06"
Set-TextValue $ws.Range("E11") "  -2.30%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.470"
Set-TextValue $ws.Range("E12") "  +5.81%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000228"
Set-TextValue $ws.Range("E13") "  +2.06%  "

# Row 14
Set-TextValue $ws.Range("D14") "35.17"
Set-TextValue $ws.Range("E14") "  +2.89%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.575.66"
Set-TextValue $ws.Range("E15") "  +1.08%  "

# Row 16
Set-TextValue $ws.Range("D16") "64.431.48"
Set-TextValue $ws.Range("E16") "  +2.89%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.076.55"
Set-TextValue $ws.Range("E17") "  +1.20%  "

# Row 18
Set-TextValue $ws.Range("E18") "  +2.12%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.77"
Set-TextValue $ws.Range("E19") "  +2.21%  "

# Row 20
Set-TextValue $ws.Range("D20") "478.22"
Set-TextValue $ws.Range("E20") "  +0.73%  "

# Row 21
Set-TextValue $ws.Range("D21") "13.93"
Set-TextValue $ws.Range("E21") "  +3.71%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.676"
Set-TextValue $ws.Range("E22") "  +1.95%  "

# Row 23
Set-TextValue $ws.Range("D23") "7.54"
Set-TextValue $ws.Range("E23") "  +6.60%  "

# Row 24
Set-TextValue $ws.Range("D24") "13.66"
Set-TextValue $ws.Range("E24") "  +11.85%  "

# Row 25
Set-TextValue $ws.Range("D25") "81.57"
Set-TextValue $ws.Range("E25") "  +2.00%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  +0.04%  "

# Row 27
Set-TextValue $ws.Range("D27") "2.80"
Set-TextValue $ws.Range("E27") "  +2.98%  "

# Row 28
Set-TextValue $ws.Range("D28") "8.10"
Set-TextValue $ws.Range("E28") "  +4.76%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.07"
Set-TextValue $ws.Range("E29") "  +6.41%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("E30") "  -0.07%  "

# Row 31
Set-TextValue $ws.Range("D31") "26.18"
Set-TextValue $ws.Range("E31") "  +2.27%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.14"
Set-TextValue $ws.Range("E32") "  +0.13%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.50"
Set-TextValue $ws.Range("E33") "  +5.73%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.57"
Set-TextValue $ws.Range("E34") "  -0.15%  "

# Row 35
Set-TextValue $ws.Range("B35") "OKB"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D35") "55.62"
Set-TextValue $ws.Range("E35") "  +1.88%  "

# Row 36
Set-TextValue $ws.Range("B36") "Filecoin"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D36") "6.18"
Set-TextValue $ws.Range("E36") "  +5.53%  "

# Row 37
Set-TextValue $ws.Range("D37") "460.31"
Set-TextValue $ws.Range("E37") "  +0.90%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.01"
Set-TextValue $ws.Range("E38") "  +21.43%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0828"
Set-TextValue $ws.Range("E39") "  +3.27%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0406"
Set-TextValue $ws.Range("E40") "  +4.38%  "

# Row 41
Set-TextValue $ws.Range("D41") "2.969.25"
Set-TextValue $ws.Range("E41") "  -3.15%  "

# Row 42
Set-TextValue $ws.Range("D42") "8.25"
Set-TextValue $ws.Range("E42") "  +1.57%  "

# Row 43
Set-TextValue $ws.Range("E43") "  -0.51%  "

# Row 44
Set-TextValue $ws.Range("D44") "27.82"
Set-TextValue $ws.Range("E44") "  +2.14%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.263"
Set-TextValue $ws.Range("E45") "  +6.60%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.16"
Set-TextValue $ws.Range("E46") "  +7.95%  "

# Row 47
Set-TextValue $ws.Range("E47") "  +0.08%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.112"
Set-TextValue $ws.Range("E48") "  +3.74%  "

# Row 49
Set-TextValue $ws.Range("D49") "120.60"
Set-TextValue $ws.Range("E49") "  +3.74%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0₃0515"
Set-TextValue $ws.Range("E50") "  +4.05%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.08"
Set-TextValue $ws.Range("E51") "  +2.34%  "
